$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overall")

# ---------------------------------------------------------------------------
# 1) Insert the new rows (blank separators + per-vendor expense-breakdown
#    rows) underneath "Home Depot" and "Contractor". Working top-down with
#    running offsets keeps every insertion point a simple constant:
#      11       -> blank row above Home Depot
#      13:16    -> 4 detail rows under Home Depot
#      17       -> blank row under the Home Depot detail rows
#      19:21    -> 3 detail rows under Contractor
#      22       -> blank row under the Contractor detail rows
# ---------------------------------------------------------------------------
$ws.Rows("11:11").Insert()
$ws.Rows("13:16").Insert()
$ws.Rows("17:17").Insert()
$ws.Rows("19:21").Insert()
$ws.Rows("22:22").Insert()

# ---------------------------------------------------------------------------
# 2) Formatting. New rows come in with no explicit style, so copy it in from
#    known-good neighbours instead of typing raw NumberFormat strings (which
#    would otherwise mint duplicate style records).
#    - Home Depot's detail/blank rows (13-17) mirror the plain bordered
#      style used throughout column A/B (same as row 10).
#    - Home Depot's own total (row 12) switches to a 2-decimal format.
#    - Contractor's total (row 18) AND all of its detail/blank rows
#      (19-22) use that same 2-decimal format.
# ---------------------------------------------------------------------------
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A13:B17").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B12").NumberFormat = "0.00"

$ws.Range("A10").Copy()
$ws.Range("A18:A22").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B12").Copy()
$ws.Range("B18:B22").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Home Depot seasonal breakdown (rows 13-16).
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = "  Water Heater"
$ws.Range("B13").Formula = "=SUM(Jan!B11, Feb!B11, Mar!B11, Apr!B11, May!B11)"

$ws.Range("A14").Value = "  Roof"
$ws.Range("B14").Formula = "=SUM(June!B11, July!B11, Aug!B11, Sep!B11)"

$ws.Range("A15").Value = "  Fence (Front, Right, Left)"
$ws.Range("B15").Formula = "=SUM(Oct!B11, Nov!B11)"

$ws.Range("A16").Value = "  Electrical Panel"
$ws.Range("B16").Formula = "=SUM(Dec!B11)"

# ---------------------------------------------------------------------------
# 4) Contractor seasonal breakdown (rows 19-21).
# ---------------------------------------------------------------------------
$ws.Range("A19").Value = "   Plumbing"
$ws.Range("B19").Formula = "=SUM(Jan!B12, Feb!B12, Mar!B12, Apr!B12, May!B12)"

$ws.Range("A20").Value = "  Roof"
$ws.Range("B20").Formula = "=SUM(June!B12, July!B12, Aug!B12, Sep!B12)"

$ws.Range("A21").Value = "  Electrical Panel"
$ws.Range("B21").Formula = "=SUM(Dec!B12)"

# ---------------------------------------------------------------------------
# 5) Column A widens to fit the longest new label ("  Fence (Front, Right,
#    Left)"). The headless engine quantizes column widths to whole
#    characters, so this is the closest achievable match to the recorded
#    24.140625 width.
# ---------------------------------------------------------------------------
$ws.Columns("A:A").ColumnWidth = 23.14

# ---------------------------------------------------------------------------
# 6) Restore the view's last active selection.
# ---------------------------------------------------------------------------
$ws.Range("E10").Select()
